$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.440895
$ws.Range("H2").Value = 73.322685
$ws.Range("I2").Value = 0.2259632546784623
$ws.Range("J2").Value = 0.2259632546784623
$ws.Range("M2").Value = 7.374758333333332
$ws.Range("N2").Value = 22.124275
$ws.Range("O2").Value = 0.2902113952021427
$ws.Range("P2").Value = 0.2902113952021427
$ws.Range("Q2").Value = 180.245694075375
$ws.Range("R2").Value = 1622.211246678375
$ws.Range("S2").Value = 0.06557711140465365
$ws.Range("T2").Value = 0.06557711140465365

$ws.Range("G3").Value = 24.440895
$ws.Range("H3").Value = 73.322685
$ws.Range("I3").Value = 0.2259632546784623
$ws.Range("J3").Value = 0.2259632546784623
$ws.Range("O3").Value = 0.07140925709849509
$ws.Range("P3").Value = 0.07140925709849508
$ws.Range("Q3").Value = 44.351156852955
$ws.Range("R3").Value = 399.1604116765951
$ws.Range("S3").Value = 0.01613586814814704
$ws.Range("T3").Value = 0.01613586814814703

$ws.Range("G4").Value = 24.440895
$ws.Range("H4").Value = 73.322685
$ws.Range("I4").Value = 0.2259632546784623
$ws.Range("J4").Value = 0.2259632546784623
$ws.Range("O4").Value = 0.6383793476993622
$ws.Range("P4").Value = 0.6383793476993621
$ws.Range("Q4").Value = 396.48728654955
$ws.Range("R4").Value = 3568.385578945951
$ws.Range("S4").Value = 0.1442502751256616
$ws.Range("T4").Value = 0.1442502751256616

$ws.Range("I5").Value = 0.6065386452756374
$ws.Range("J5").Value = 0.6065386452756373
$ws.Range("M5").Value = 7.374758333333332
$ws.Range("N5").Value = 22.124275
$ws.Range("O5").Value = 0.2902113952021427
$ws.Range("P5").Value = 0.2902113952021427
$ws.Range("Q5").Value = 483.8219349283667
$ws.Range("R5").Value = 4354.3974143553
$ws.Range("S5").Value = 0.1760244264894603
$ws.Range("T5").Value = 0.1760244264894602

$ws.Range("I6").Value = 0.6065386452756374
$ws.Range("J6").Value = 0.6065386452756373
$ws.Range("O6").Value = 0.07140925709849509
$ws.Range("P6").Value = 0.07140925709849508
$ws.Range("S6").Value = 0.0433124740606609
$ws.Range("T6").Value = 0.04331247406066089

$ws.Range("I7").Value = 0.6065386452756374
$ws.Range("J7").Value = 0.6065386452756373
$ws.Range("O7").Value = 0.6383793476993622
$ws.Range("P7").Value = 0.6383793476993621
$ws.Range("S7").Value = 0.3872017447255163
$ws.Range("T7").Value = 0.3872017447255161

$ws.Range("I8").Value = 0.1674981000459004
$ws.Range("J8").Value = 0.1674981000459004
$ws.Range("M8").Value = 7.374758333333332
$ws.Range("N8").Value = 22.124275
$ws.Range("O8").Value = 0.2902113952021427
$ws.Range("P8").Value = 0.2902113952021427
$ws.Range("Q8").Value = 133.6093841542528
$ws.Range("R8").Value = 1202.484457388275
$ws.Range("S8").Value = 0.04860985730802884
$ws.Range("T8").Value = 0.04860985730802884

$ws.Range("I9").Value = 0.1674981000459004
$ws.Range("J9").Value = 0.1674981000459004
$ws.Range("O9").Value = 0.07140925709849509
$ws.Range("P9").Value = 0.07140925709849508
$ws.Range("S9").Value = 0.01196091488968716
$ws.Range("T9").Value = 0.01196091488968715

$ws.Range("I10").Value = 0.1674981000459004
$ws.Range("J10").Value = 0.1674981000459004
$ws.Range("O10").Value = 0.6383793476993622
$ws.Range("P10").Value = 0.6383793476993621
$ws.Range("S10").Value = 0.1069273278481844
$ws.Range("T10").Value = 0.1069273278481844
